# The "Result Footer: Displayed When The User HAS NOT Voted" template (row 9,
# column C) previously carried its own long HTML fragment that duplicated the
# "Total votes" footer but also re-added a "Vote" button/paragraph. The
# commit collapses it to the same short footer markup already used in row 8
# (Result Footer: Displayed When The User HAS Voted), so both footers share
# one template string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value2 = $ws.Range("C8").Value2

# Leave the selection on the cell that was edited, like Excel would after an
# interactive edit.
$ws.Range("C9").Select() | Out-Null
